$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 25114
$ws.Range("B3").Value = 526
$ws.Range("B4").Value = 506
$ws.Range("B5").Value = 3428
$ws.Range("B6").Value = 2248
$ws.Range("B7").Value = 1173
$ws.Range("B8").Value = 138
$ws.Range("B9").Value = 899
$ws.Range("B10").Value = 423
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 86
$ws.Range("B13").Value = 208
